# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 281
$ws1.Range("F5").Value = 156
$ws1.Range("F6").Value = 94
$ws1.Range("F7").Value = 275
$ws1.Range("F9").Value = 2012
$ws1.Range("F10").Value = 354
$ws1.Range("F11").Value = 4771
$ws1.Range("F13").Value = 336

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 281
$ws4.Range("F7").Value = 156
$ws4.Range("F8").Value = 94
$ws4.Range("F9").Value = 275
$ws4.Range("F13").Value = 2012
$ws4.Range("F14").Value = 354
$ws4.Range("F15").Value = 4771
$ws4.Range("F17").Value = 336

$wb.Save()
